$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (hunk 0)
$ws.Range("H28").Value = 454
$ws.Range("J28").Value = 896.6667
$ws.Range("L28").Value = 896.6667
$ws.Range("N28").Value = -1866.6667

# Row 88 (hunk 1)
$ws.Range("H88").Value = 8334732
$ws.Range("I88").Value = 14286799
$ws.Range("K88").Value = 14286799
$ws.Range("M88").Value = -14286393

# Row 91 (hunk 2)
$ws.Range("H91").Value = 8334732
$ws.Range("I91").Value = 14286799
$ws.Range("K91").Value = 14286799
$ws.Range("M91").Value = -14285395

# Row 98 (hunk 3)
$ws.Range("H98").Value = 1775
$ws.Range("I98").Value = 1880.3125
$ws.Range("J98").Value = 90
$ws.Range("K98").Value = 1880.3125
$ws.Range("L98").Value = 90
$ws.Range("M98").Value = -382.3125
$ws.Range("N98").Value = -3086

# Row 122 (hunk 4)
$ws.Range("H122").Value = 1775
$ws.Range("I122").Value = 1880.3125
$ws.Range("J122").Value = 90
$ws.Range("K122").Value = 5640.9375
$ws.Range("L122").Value = 270
$ws.Range("M122").Value = -3190.9375
$ws.Range("N122").Value = -5170

# Row 129 (hunk 5)
$ws.Range("H129").Value = 2271.5
$ws.Range("I129").Value = 1480.2
$ws.Range("K129").Value = 4440.6
$ws.Range("M129").Value = 559.3999999999996

# Row 137 (hunk 6)
$ws.Range("H137").Value = 2389113
$ws.Range("I137").Value = 4546811
$ws.Range("K137").Value = 13640433
$ws.Range("M137").Value = -13637883

# Row 138 (hunk 7)
$ws.Range("H138").Value = 6087.737
$ws.Range("I138").Value = 2652.077
$ws.Range("J138").Value = 7102.8184
$ws.Range("K138").Value = 7956.231000000001
$ws.Range("L138").Value = 21308.4552
$ws.Range("M138").Value = -2816.231000000001
$ws.Range("N138").Value = -31588.4552

$ws = $wb.Worksheets.Item("ARM")
# Row 97 (hunk 8)
$ws.Range("H97").Value = 3026.2307
$ws.Range("I97").Value = 3144.7273
$ws.Range("J97").Value = 2374.5
$ws.Range("K97").Value = 3144.7273
$ws.Range("L97").Value = 2374.5
$ws.Range("M97").Value = -2648.7273
$ws.Range("N97").Value = -3366.5

# Row 122 (hunk 9)
$ws.Range("H122").Value = 3970.842
$ws.Range("I122").Value = 3840.7273
$ws.Range("J122").Value = 4149.75
$ws.Range("K122").Value = 11522.1819
$ws.Range("L122").Value = 12449.25
$ws.Range("M122").Value = -9072.1819
$ws.Range("N122").Value = -17349.25

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (hunk 10)
$ws.Range("H99").Value = 3343.6365
$ws.Range("I99").Value = 2928
$ws.Range("K99").Value = 2928
$ws.Range("M99").Value = -1430

# Row 105 (hunk 11)
$ws.Range("H105").Value = 10834809
$ws.Range("I105").Value = 527285
$ws.Range("K105").Value = 527285
$ws.Range("M105").Value = -525538

# Row 134 (hunk 12)
$ws.Range("H134").Value = 2854.0625
$ws.Range("I134").Value = 2194
$ws.Range("J134").Value = 3154.0908
$ws.Range("K134").Value = 6582
$ws.Range("L134").Value = 9462.2724
$ws.Range("M134").Value = -4047
$ws.Range("N134").Value = -14532.2724

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (hunk 13)
$ws.Range("H31").Value = 2849.047
$ws.Range("I31").Value = 1403.2778
$ws.Range("J31").Value = 3237.4626
$ws.Range("K31").Value = 1403.2778
$ws.Range("L31").Value = 3237.4626
$ws.Range("M31").Value = -1108.2778
$ws.Range("N31").Value = -3827.4626

# Row 34 (hunk 14)
$ws.Range("H34").Value = 2849.047
$ws.Range("I34").Value = 1403.2778
$ws.Range("J34").Value = 3237.4626
$ws.Range("K34").Value = 1403.2778
$ws.Range("L34").Value = 3237.4626
$ws.Range("M34").Value = -1201.2778
$ws.Range("N34").Value = -3641.4626

# Row 58 (hunk 15)
$ws.Range("H58").Value = 5127.6
$ws.Range("I58").Value = 4550.1665
$ws.Range("K58").Value = 4550.1665
$ws.Range("M58").Value = -4347.1665

# Row 99 (hunk 16)
$ws.Range("H99").Value = 76926350
$ws.Range("I99").Value = 90911700
$ws.Range("K99").Value = 90911700
$ws.Range("M99").Value = -90910202

# Row 107 (hunk 17)
$ws.Range("H107").Value = 725.76666
$ws.Range("I107").Value = 356.13635
$ws.Range("J107").Value = 1742.25
$ws.Range("K107").Value = 356.13635
$ws.Range("L107").Value = 1742.25
$ws.Range("M107").Value = 1563.86365
$ws.Range("N107").Value = -5582.25

# Row 126 (hunk 18)
$ws.Range("H126").Value = 76926350
$ws.Range("I126").Value = 90911700
$ws.Range("K126").Value = 272735100
$ws.Range("M126").Value = -272732630

# Row 136 (hunk 19)
$ws.Range("H136").Value = 5127.6
$ws.Range("I136").Value = 4550.1665
$ws.Range("K136").Value = 13650.4995
$ws.Range("M136").Value = -11100.4995

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (hunk 20)
$ws.Range("H4").Value = 1956591.4
$ws.Range("I4").Value = 1517319.8
$ws.Range("K4").Value = 4551959.4
$ws.Range("M4").Value = -4551847.4

# Row 12 (hunk 21)
$ws.Range("H12").Value = 56.933334
$ws.Range("I12").Value = 5.5
$ws.Range("J12").Value = 91.22221999999999
$ws.Range("K12").Value = 16.5
$ws.Range("L12").Value = 273.66666
$ws.Range("M12").Value = 156.5
$ws.Range("N12").Value = -619.66666

# Row 131 (hunk 22)
$ws.Range("H131").Value = 8398
$ws.Range("I131").Value = 26248.625
$ws.Range("K131").Value = 78745.875
$ws.Range("M131").Value = -73705.875

# Row 132 (hunk 23)
$ws.Range("H132").Value = 5477.881
$ws.Range("I132").Value = 5986.625
$ws.Range("K132").Value = 53879.625
$ws.Range("M132").Value = -51349.625

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (hunk 24)
$ws.Range("H2").Value = 37
$ws.Range("I2").Value = 43
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 43
$ws.Range("L2").Value = 16
$ws.Range("M2").Value = 70
$ws.Range("N2").Value = -242

# Row 58 (hunk 25)
$ws.Range("H58").Value = 42000
$ws.Range("J58").Value = 42000
$ws.Range("L58").Value = 42000
$ws.Range("N58").Value = -42554

# Row 97 (hunk 26)
$ws.Range("H97").Value = 721.8148
$ws.Range("I97").Value = 707.17645
$ws.Range("J97").Value = 746.7
$ws.Range("K97").Value = 707.17645
$ws.Range("L97").Value = 746.7
$ws.Range("M97").Value = -211.17645
$ws.Range("N97").Value = -1738.7

# Row 122 (hunk 27)
$ws.Range("H122").Value = 105268264
$ws.Range("I122").Value = 142859340
$ws.Range("K122").Value = 428578020
$ws.Range("M122").Value = -428575570

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (hunk 28)
$ws.Range("H46").Value = 2000
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -3376

# Row 95 (hunk 29)
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("M95").ClearContents()
$ws.Range("N95").ClearContents()

# Row 100 (hunk 30)
$ws.Range("H100").Value = 4000
$ws.Range("I100").Value = 4000
$ws.Range("K100").Value = 4000
$ws.Range("M100").Value = -3459

# Row 122 (hunk 31)
$ws.Range("H122").Value = 3800
$ws.Range("I122").Value = 3800
$ws.Range("K122").Value = 11400
$ws.Range("M122").Value = -8950

# Row 132 (hunk 32)
$ws.Range("H132").Value = 3958.1562
$ws.Range("I132").Value = 3256.2307
$ws.Range("J132").Value = 6999.8335
$ws.Range("K132").Value = 9768.6921
$ws.Range("L132").Value = 20999.5005
$ws.Range("M132").Value = -7238.6921
$ws.Range("N132").Value = -26059.5005

# Row 136 (hunk 33)
$ws.Range("H136").Value = 7134.826
$ws.Range("I136").Value = 5147.375
$ws.Range("K136").Value = 15442.125
$ws.Range("M136").Value = -12892.125

# Row 137 (hunk 34)
$ws.Range("H137").Value = 59089.91
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# Row 139 (hunk 35)
$ws.Range("H139").Value = 70075
$ws.Range("J139").Value = 70075
$ws.Range("L139").Value = 70075
$ws.Range("N139").Value = -80355

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (hunk 36)
$ws.Range("H96").Value = 1819.3636
$ws.Range("I96").Value = 1776.75
$ws.Range("J96").Value = 1933
$ws.Range("K96").Value = 1776.75
$ws.Range("L96").Value = 1933
$ws.Range("M96").Value = -403.75
$ws.Range("N96").Value = -4679

# Row 100 (hunk 37)
$ws.Range("H100").Value = 31250968
$ws.Range("I100").Value = 687.52
$ws.Range("J100").Value = 142859120
$ws.Range("K100").Value = 1375.04
$ws.Range("L100").Value = 285718240
$ws.Range("M100").Value = -834.04
$ws.Range("N100").Value = -285719322
